{"js": "// Remove the stray \"_GoBack\" bookmark that Word drops at the cursor's last\n// edit position. Deleting it also makes Word renumber the remaining\n// bookmark ids so they stay contiguous, matching a normal re-save.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the stray \"_GoBack\" bookmark left over from the cursor's last edit\n# position. Word auto-renumbers the remaining bookmark ids to stay\n# contiguous, matching a normal re-save in a newer Word build.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
